$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 to hold the "water" record (previously row 7's data),
# now the only remaining data row.
$ws.Range("A2").Value = 7
$ws.Range("B2").Value = "water"
$ws.Range("C2").Value = 34
$ws.Range("D2").Value = "D:\work-place\flutter apps\projects\sharp\BisleriumCafeBackend\fyp-document\fyp\coffee\coffee-image\2024-01-08\1704735419442-5ce111d5-dd7c-41f3-b432-7abff9a14dd6.jpg"

# Remove the now-obsolete rows 3 through 7.
$ws.Range("A3:D7").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftUp)

# Update the active selection to match the saved view state.
$ws.Range("H14").Select()
